$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '34.738.47'
$ws.Range("E2").Value = '  -0.48%  '
$ws.Range("D3").Value = '1.819.67'
$ws.Range("E3").Value = '  +0.45%  '
$ws.Range("E4").Value = '  +0.62%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '230.32'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.78%  '
$ws.Range("E6").Value = '  +0.79%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '39.35'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -2.73%  '
$ws.Range("E9").Value = '  +1.47%  '
$ws.Range("E10").Value = '  -0.50%  '
$ws.Range("E11").Value = '  -0.77%  '
$ws.Range("D12").Value = '2.084.11'
$ws.Range("E12").Value = '  +0.57%  '
$ws.Range("B13").Value = 'WrappedEther'
$ws.Range("C13").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D13").Value = '1.831.48'
$ws.Range("E13").Value = '  +1.09%  '
$ws.Range("B14").Value = 'Chainlink'
$ws.Range("C14").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '11.24'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +1.22%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.663'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +0.96%  '
$ws.Range("E16").Value = '  -1.60%  '
$ws.Range("D17").Value = '34.618.04'
$ws.Range("E17").Value = '  -0.67%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '69.35'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +0.32%  '
$ws.Range("D19").Value = '0.0₃0782'
$ws.Range("E19").Value = '  -0.27%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '238.61'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +0.77%  '
$ws.Range("E21").Value = '  +1.64%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '4.62'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -0.45%  '
$ws.Range("E23").Value = '  +0.46%  '
$ws.Range("E24").Value = '  -0.68%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '172.89'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +0.18%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '7.71'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -2.26%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.124'
$ws.Range("D27").Style = "Normal"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '17.25'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -0.97%  '
$ws.Range("E29").Value = '  -7.04%  '
$ws.Range("E30").Value = '  +0.63%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.0546'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -1.09%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.89'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +0.28%  '
$ws.Range("E33").Value = '  -1.01%  '
$ws.Range("E34").Value = '  +7.86%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.82'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +1.92%  '
$ws.Range("E36").Value = '  +12.05%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.694'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +2.31%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '91.01'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -2.44%  '
$ws.Range("D39").Value = '1.338.10'
$ws.Range("E39").Value = '  +2.38%  '
$ws.Range("E40").Value = '  +2.19%  '
$ws.Range("E41").Value = '  -0.11%  '
$ws.Range("E42").Value = '  -3.10%  '
$ws.Range("E43").Value = '  -0.49%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '2.23'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -4.18%  '
$ws.Range("E45").Value = '  -0.39%  '
$ws.Range("B46").Value = 'Kaspa'
$ws.Range("C46").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.0522'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +1.84%  '
$ws.Range("B47").Value = 'FraxShare'
$ws.Range("C47").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '6.24'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -0.22%  '
$ws.Range("D48").Value = '1.999.96'
$ws.Range("E48").Value = '  +0.63%  '
$ws.Range("E49").Value = '  +0.56%  '
$ws.Range("E50").Value = '  +3.83%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '97.89'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -1.33%  '

Write-Host "Applied all cell updates"